# Weekly refresh of the "Hortaliza, Mapocho Venta Directa de Santiago - Sandia"
# sheet: the 21 existing price records (rows 2-22) are reshuffled into new
# row positions (same records, same values, just relocated) as part of the
# weekly consolidation pass. No values are altered - only which physical
# sheet row each record occupies changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of destination row -> source row (1-based data rows, header is row 1).
$rowMap = @{
    2  = 22
    3  = 16
    4  = 10
    5  = 9
    6  = 5
    7  = 6
    8  = 7
    9  = 8
    10 = 20
    11 = 13
    12 = 14
    13 = 2
    14 = 3
    15 = 15
    16 = 12
    17 = 4
    18 = 17
    19 = 18
    20 = 19
    21 = 11
    22 = 21
}

$firstCol = 1   # A
$lastCol  = 18  # R

# Snapshot every source row's values before any writes happen, so that
# overwriting one row never clobbers data still needed for another.
$snapshot = @{}
for ($r = 2; $r -le 22; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c - $firstCol]
    }
}
